# Update the "Equal Rows" / "Different Rows" summary counts (and their
# relative/percentage counterparts) on the full, left, right and inner
# comparison-summary sheets to reflect the newly added column-limit tests.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "full";  B4 = 395.0; C2 = 0.6046046046046046;  C4 = 0.3953953953953954 },
    @{ Sheet = "left";  B4 = 363.0; C2 = 0.6246122026887281;  C4 = 0.375387797311272 },
    @{ Sheet = "right"; B4 = 388.0; C2 = 0.6088709677419355;  C4 = 0.3911290322580645 },
    @{ Sheet = "inner"; B4 = 356.0; C2 = 0.6291666666666667;  C4 = 0.37083333333333335 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    $ws.Range("B2").Value = 604.0
    $ws.Range("C2").Value = $u.C2

    $ws.Range("B4").Value = $u.B4
    $ws.Range("C4").Value = $u.C4
}
